$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.081.63"
$ws.Range("E2").Value = "  +7.30%  "
$ws.Range("D3").Value = "2.576.86"
$ws.Range("E3").Value = "  +9.62%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "504.24"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +5.86%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "156.25"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +7.60%  "
$ws.Range("E7").Value = "  +3.98%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.993"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").Value = "2.570.83"
$ws.Range("E9").Value = "  +9.08%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.15"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +13.34%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.103"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.86%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.340"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +4.96%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.127"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "2.972.97"
$ws.Range("E14").Value = "  +7.71%  "
$ws.Range("D15").Value = "59.010.71"
$ws.Range("E15").Value = "  +7.37%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.84"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +8.65%  "
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").Value = "2.563.86"
$ws.Range("E18").Value = "  +9.12%  "
$ws.Range("E19").Value = "  +3.55%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "334.19"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +5.37%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.32"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.85%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.04"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +6.95%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "59.65"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +5.03%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.416"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +5.37%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.166"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +7.43%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.992"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.646.32"
$ws.Range("E28").Value = "  +8.23%  "
$ws.Range("D29").Value = "0.0₃0828"
$ws.Range("E29").Value = "  +9.10%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.87%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "157.55"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +7.02%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.37"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.77%  "
$ws.Range("E34").Value = "  +5.92%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.51"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +8.61%  "
$ws.Range("E36").Value = "  +9.41%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.90"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +8.60%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.848"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("E39").Value = "  +10.54%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.44"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.21%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "35.09"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.37%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "291.07"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +15.54%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.102"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.94%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.624"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +8.23%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0566"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +7.93%  "
$ws.Range("B46").Value = "SuiNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.776"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +22.69%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "19.25"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +14.35%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "4.85"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +9.42%  "
$ws.Range("E50").Value = "  +6.10%  "
$ws.Range("D51").Value = "1.997.94"
$ws.Range("E51").Value = "  +12.05%  "
